# Updates cryptos list price/volume figures (and two re-ordered coin rows)
# to match the Fri Nov 10 08:40:50 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.490.37"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "2.095.25"
$ws.Range("E3").Value = "  +9.53%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'251.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").Value = "'0.654"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.58%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'47.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.25%  "
$ws.Range("D9").Value = "'59.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("D10").Value = "'0.372"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").Value = "'0.0742"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("D12").Value = "'0.0997"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'14.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.397.57"
$ws.Range("E14").Value = "  +9.40%  "
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "2.094.01"
$ws.Range("E16").Value = "  +9.19%  "
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "36.467.69"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").Value = "'72.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.56%  "
$ws.Range("D20").Value = "0.0₃0826"
$ws.Range("E20").Value = "  -4.17%  "
$ws.Range("D21").Value = "'13.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("D22").Value = "'238.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("D23").Value = "'5.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "'2.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.77%  "
$ws.Range("D26").Value = "'169.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("D27").Value = "'21.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.39%  "
$ws.Range("D28").Value = "'9.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.12%  "
$ws.Range("E29").Value = "  -10.34%  "
$ws.Range("D30").Value = "'27.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +57.24%  "
$ws.Range("E31").Value = "  -5.25%  "
$ws.Range("D32").Value = "'4.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").Value = "'0.0607"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").Value = "'0.0926"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.82%  "
$ws.Range("D35").Value = "'0.984"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.66%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +15.52%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'1.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("E39").Value = "  -6.62%  "
$ws.Range("E40").Value = "  -12.21%  "
$ws.Range("E41").Value = "  +4.70%  "
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("D43").Value = "'96.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.70%  "
$ws.Range("E44").Value = "  -6.59%  "
$ws.Range("D45").Value = "'15.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.32%  "
$ws.Range("D46").Value = "1.325.99"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").Value = "'0.0837"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.10%  "
$ws.Range("D48").Value = "'6.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.45%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.290.45"
$ws.Range("E49").Value = "  +9.32%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("D51").Value = "'2.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.65%  "
